# Remove the parentheses around the area code in every phone number
# (column F), e.g. "(291) 553-0508" -> "291 553-0508".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val -ne $null -and $val -match '^\((\d+)\)\s*(.*)$') {
        $cell.Value = "$($Matches[1]) $($Matches[2])"
    }
}

# The last three guest rows no longer need the extra wrapped-text height
# once their phone numbers shrank - let those rows re-fit to the default.
foreach ($r in 9, 10, 11) {
    if ($r -le $lastRow) {
        $ws.Rows.Item($r).AutoFit()
    }
}

# Update the active selection to match the authored state.
$ws.Range("K4").Select()
